$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Proirity" -> "Priority" header typo (cell E6 in the shared-strings table).
$ws.Range("E6").Value = "Priority"

# Move the saved selection from E17 to E11.
[void]$ws.Range("E11").Select()
